$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Iceland Premier League")

# Row 45 <- old row 46 data
$ws.Range("B45").Value2 = 6102213
$ws.Range("E45").Value2 = "KR Reykjavik"
$ws.Range("F45").Value2 = "Valur Reykjavik"
$ws.Range("G45").Value2 = 0
$ws.Range("H45").Value2 = 4
$ws.Range("J45").Value2 = 2.875
$ws.Range("K45").Value2 = 3.5
$ws.Range("L45").Value2 = 2.3
$ws.Range("M45").Value2 = 2.7
$ws.Range("N45").Value2 = 3.6
$ws.Range("O45").Value2 = 2.45
$ws.Range("P45").Value2 = 0
$ws.Range("Q45").Value2 = 2
$ws.Range("T45").Value2 = 1.875
$ws.Range("U45").Value2 = 1.975
$ws.Range("X45").Value2 = 1.45
$ws.Range("AA45").Value2 = 0.875

# Row 46 <- old row 45 data
$ws.Range("B46").Value2 = 6102350
$ws.Range("E46").Value2 = "Keflavik"
$ws.Range("F46").Value2 = "FH Hafnarfjordur"
$ws.Range("G46").Value2 = 2
$ws.Range("H46").Value2 = 3
$ws.Range("J46").Value2 = 3.8
$ws.Range("K46").Value2 = 3.6
$ws.Range("L46").Value2 = 1.909
$ws.Range("M46").Value2 = 3.8
$ws.Range("N46").Value2 = 3.75
$ws.Range("O46").Value2 = 1.85
$ws.Range("P46").Value2 = 0.5
$ws.Range("Q46").Value2 = 1.95
$ws.Range("T46").Value2 = 1.95
$ws.Range("U46").Value2 = 1.85
$ws.Range("X46").Value2 = 0.8500000000000001
$ws.Range("AA46").Value2 = 0.95

# Row 91 <- old row 94 data
$ws.Range("B91").Value2 = 7173166
$ws.Range("E91").Value2 = "Valur Reykjavik"
$ws.Range("F91").Value2 = "Breidablik"
$ws.Range("G91").Value2 = 4
$ws.Range("H91").Value2 = 2
$ws.Range("J91").Value2 = 2.25
$ws.Range("K91").Value2 = 3.6
$ws.Range("L91").Value2 = 2.6
$ws.Range("M91").Value2 = 2.7
$ws.Range("N91").Value2 = 4
$ws.Range("O91").Value2 = 2.1
$ws.Range("P91").Value2 = 0.25
$ws.Range("Q91").Value2 = 1.925
$ws.Range("R91").Value2 = 1.875
$ws.Range("S91").Value2 = 3.75
$ws.Range("T91").Value2 = 1.925
$ws.Range("U91").Value2 = 1.775
$ws.Range("V91").Value2 = 1.7
$ws.Range("Y91").Value2 = 0.925
$ws.Range("Z91").Value2 = -1
$ws.Range("AA91").Value2 = 0.925
$ws.Range("AB91").Value2 = -1

# Row 92 <- old row 93 data
$ws.Range("B92").Value2 = 7173183
$ws.Range("E92").Value2 = "HK Kopavogur"
$ws.Range("F92").Value2 = "Fylkir Reykjavik"
$ws.Range("H92").Value2 = 2
$ws.Range("I92").Value2 = "D"
$ws.Range("J92").Value2 = 2.2
$ws.Range("K92").Value2 = 3.6
$ws.Range("L92").Value2 = 2.7
$ws.Range("M92").Value2 = 2.45
$ws.Range("N92").Value2 = 3.6
$ws.Range("O92").Value2 = 2.4
$ws.Range("P92").Value2 = 0
$ws.Range("Q92").Value2 = 1.925
$ws.Range("R92").Value2 = 1.925
$ws.Range("T92").Value2 = 1.9
$ws.Range("U92").Value2 = 1.95
$ws.Range("V92").Value2 = -1
$ws.Range("W92").Value2 = 2.6
$ws.Range("Y92").Value2 = 0
$ws.Range("Z92").Value2 = 0
$ws.Range("AA92").Value2 = 0.8999999999999999
$ws.Range("AB92").Value2 = -1

# Row 93 <- old row 95 data
$ws.Range("B93").Value2 = 7173182
$ws.Range("E93").Value2 = "Fram Reykjavik"
$ws.Range("F93").Value2 = "Keflavik"
$ws.Range("G93").Value2 = 3
$ws.Range("H93").Value2 = 1
$ws.Range("I93").Value2 = "H"
$ws.Range("J93").Value2 = 1.75
$ws.Range("K93").Value2 = 3.8
$ws.Range("L93").Value2 = 3.75
$ws.Range("M93").Value2 = 1.833
$ws.Range("N93").Value2 = 3.8
$ws.Range("O93").Value2 = 3.5
$ws.Range("P93").Value2 = -0.5
$ws.Range("Q93").Value2 = 1.825
$ws.Range("R93").Value2 = 1.975
$ws.Range("T93").Value2 = 1.925
$ws.Range("U93").Value2 = 1.875
$ws.Range("V93").Value2 = 0.833
$ws.Range("W93").Value2 = -1
$ws.Range("Y93").Value2 = 0.825
$ws.Range("Z93").Value2 = -1
$ws.Range("AA93").Value2 = 0.925

# Row 94 <- old row 91 data
$ws.Range("B94").Value2 = 7173167
$ws.Range("E94").Value2 = "Vikingur Reykjavik"
$ws.Range("F94").Value2 = "FH Hafnarfjordur"
$ws.Range("G94").Value2 = 2
$ws.Range("H94").Value2 = 1
$ws.Range("J94").Value2 = 1.6
$ws.Range("K94").Value2 = 4
$ws.Range("L94").Value2 = 4.333
$ws.Range("M94").Value2 = 1.571
$ws.Range("N94").Value2 = 4.2
$ws.Range("O94").Value2 = 4.5
$ws.Range("P94").Value2 = -1
$ws.Range("Q94").Value2 = 1.875
$ws.Range("R94").Value2 = 1.925
$ws.Range("S94").Value2 = 3.25
$ws.Range("T94").Value2 = 1.775
$ws.Range("U94").Value2 = 1.925
$ws.Range("V94").Value2 = 0.571
$ws.Range("Y94").Value2 = 0
$ws.Range("Z94").Value2 = 0
$ws.Range("AA94").Value2 = -0.5
$ws.Range("AB94").Value2 = 0.4625

# Row 95 <- old row 92 data
$ws.Range("B95").Value2 = 7173168
$ws.Range("E95").Value2 = "Stjarnan"
$ws.Range("F95").Value2 = "KR Reykjavik"
$ws.Range("G95").Value2 = 2
$ws.Range("H95").Value2 = 0
$ws.Range("J95").Value2 = 1.95
$ws.Range("K95").Value2 = 3.75
$ws.Range("L95").Value2 = 3.2
$ws.Range("M95").Value2 = 1.85
$ws.Range("O95").Value2 = 3.4
$ws.Range("T95").Value2 = 1.825
$ws.Range("U95").Value2 = 1.975
$ws.Range("V95").Value2 = 0.8500000000000001
$ws.Range("AA95").Value2 = -1
$ws.Range("AB95").Value2 = 0.9750000000000001

